# Update the "Förändrad" (changed) date column (C) for rows 2-6
# from 2023-10-22 (45221) to 2023-10-25 (45224).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C6").Value = 45224
